$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($sheet, $addr, $val) {
    $cell = $sheet.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-CellText $ws "D2" '44.957.73'
Set-CellText $ws "E2" '  +2.15%  '
Set-CellText $ws "D3" '2.359.90'
Set-CellText $ws "E3" '  +0.74%  '
Set-CellText $ws "E4" '  -0.33%  '
Set-CellText $ws "D5" '311.64'
Set-CellText $ws "E5" '  -0.16%  '
Set-CellText $ws "D6" '107.42'
Set-CellText $ws "E6" '  -0.72%  '
Set-CellText $ws "E7" '  +0.27%  '
Set-CellText $ws "E8" '  -0.22%  '
Set-CellText $ws "D9" '0.606'
Set-CellText $ws "E9" '  -2.32%  '
Set-CellText $ws "D10" '40.65'
Set-CellText $ws "E10" '  -1.94%  '
Set-CellText $ws "D11" '0.0912'
Set-CellText $ws "E11" '  -0.81%  '
Set-CellText $ws "D12" '8.40'
Set-CellText $ws "E12" '  -1.89%  '
Set-CellText $ws "E13" '  +1.23%  '
Set-CellText $ws "D14" '0.970'
Set-CellText $ws "E14" '  -3.70%  '
Set-CellText $ws "D15" '2.714.06'
Set-CellText $ws "E15" '  +0.59%  '
Set-CellText $ws "E16" '  -2.26%  '
Set-CellText $ws "D17" '2.353.79'
Set-CellText $ws "E17" '  +1.00%  '
Set-CellText $ws "D18" '44.853.90'
Set-CellText $ws "E18" '  +2.11%  '
Set-CellText $ws "D19" '14.24'
Set-CellText $ws "E19" '  +10.05%  '
Set-CellText $ws "D20" '7.19'
Set-CellText $ws "E20" '  -4.83%  '
Set-CellText $ws "E21" '  -1.54%  '
Set-CellText $ws "D22" '72.82'
Set-CellText $ws "E22" '  -1.85%  '
Set-CellText $ws "D23" '3.50'
Set-CellText $ws "E23" '  +1.35%  '
Set-CellText $ws "D24" '257.49'
Set-CellText $ws "E24" '  -4.15%  '
Set-CellText $ws "D25" '2.30'
Set-CellText $ws "E25" '  +0.78%  '
Set-CellText $ws "E26" '  -0.04%  '
Set-CellText $ws "E27" '  -0.81%  '
Set-CellText $ws "E28" '  -5.62%  '
Set-CellText $ws "E29" '  +1.28%  '
Set-CellText $ws "D30" '22.26'
Set-CellText $ws "E30" '  -1.67%  '
Set-CellText $ws "D31" '0.0958'
Set-CellText $ws "E31" '  +7.78%  '
Set-CellText $ws "D32" '37.11'
Set-CellText $ws "E32" '  -5.21%  '
Set-CellText $ws "D33" '167.56'
Set-CellText $ws "E33" '  -0.63%  '
Set-CellText $ws "E34" '  +3.90%  '
Set-CellText $ws "D35" '0.130'
Set-CellText $ws "E35" '  -2.03%  '
Set-CellText $ws "E36" '  +0.74%  '
Set-CellText $ws "D37" '4.68'
Set-CellText $ws "E37" '  -1.38%  '
Set-CellText $ws "D38" '3.91'
Set-CellText $ws "E38" '  +3.40%  '
Set-CellText $ws "E39" '  -0.90%  '
Set-CellText $ws "E40" '  -3.41%  '
Set-CellText $ws "E41" '  +1.29%  '
Set-CellText $ws "D42" '99.58'
Set-CellText $ws "E42" '  -4.60%  '
Set-CellText $ws "D43" '69.41'
Set-CellText $ws "E43" '  -3.11%  '
Set-CellText $ws "D44" '1.866.89'
Set-CellText $ws "E44" '  +11.56%  '
Set-CellText $ws "E45" '  -4.99%  '
Set-CellText $ws "E46" '  -0.34%  '
Set-CellText $ws "D47" '12.75'
Set-CellText $ws "E47" '  -4.99%  '
Set-CellText $ws "B48" 'ordi'
Set-CellText $ws "C48" 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
Set-CellText $ws "D48" '80.78'
Set-CellText $ws "E48" '  +4.98%  '
Set-CellText $ws "B49" 'THORChain'
Set-CellText $ws "C49" 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
Set-CellText $ws "D49" '5.60'
Set-CellText $ws "E49" '  +8.14%  '
Set-CellText $ws "D50" '110.50'
Set-CellText $ws "E50" '  -3.34%  '
Set-CellText $ws "E51" '  +1.83%  '
